$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.855.15"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.689.80"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'647.23"
$ws.Range("D6").Value = "'161.68"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "'0.443"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "4.311.44"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'32.76"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "3.695.09"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "69.843.09"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'16.00"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "'10.36"
$ws.Range("E20").Value = "  +5.47%  "
$ws.Range("D21").Value = "'471.54"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'80.07"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("D24").Value = "3.836.96"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000127"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "'9.19"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").Value = "'1.72"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "'0.168"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "'26.75"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "3.686.10"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'8.45"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "'180.39"
$ws.Range("E39").Value = "  +7.71%  "
$ws.Range("E40").Value = "  -5.41%  "
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  -1.20%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'29.23"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'46.66"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("D49").Value = "'1.26"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'1.06"
$ws.Range("E51").Value = "  -3.19%  "
